{"js": "// Re-save the document the way a newer Word/OOXML writer (POI 4.0.1) would:\n// the content/values are unchanged, only the low-level attribute order used\n// when Word rewrites the XML for the parts it touches. We \"touch\" every part\n// that the target writer re-serialized by nudging properties back to their\n// own current value (no semantic change), which forces Word to rewrite those\n// elements with the newer attribute ordering.\n\n// --- word/document.xml : re-emit the page size / margins (w:sectPr) ---\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst pageSetup = section.pageSetup;\npageSetup.load([\n  \"pageWidth\",\n  \"pageHeight\",\n  \"topMargin\",\n  \"rightMargin\",\n  \"bottomMargin\",\n  \"leftMargin\",\n  \"headerDistance\",\n  \"footerDistance\",\n  \"gutter\"\n]);\nawait context.sync();\n\npageSetup.pageWidth = pageSetup.pageWidth;\npageSetup.pageHeight = pageSetup.pageHeight;\npageSetup.topMargin = pageSetup.topMargin;\npageSetup.rightMargin = pageSetup.rightMargin;\npageSetup.bottomMargin = pageSetup.bottomMargin;\npageSetup.leftMargin = pageSetup.leftMargin;\npageSetup.headerDistance = pageSetup.headerDistance;\npageSetup.footerDistance = pageSetup.footerDistance;\npageSetup.gutter = pageSetup.gutter;\nawait context.sync();\n\n// --- word/styles.xml : re-emit docDefaults / latentStyles / every w:style ---\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < styles.items.length; i++) {\n  styles.items[i].load(\"priority\");\n}\nawait context.sync();\n\nfor (let i = 0; i < styles.items.length; i++) {\n  const s = styles.items[i];\n  s.priority = s.priority;\n}\nawait context.sync();\n", "ps1": "# Re-save the document the way a newer Word/OOXML writer (POI 4.0.1) would:\n# the content/values are unchanged, only the low-level attribute order used\n# when Word rewrites the XML for the parts it touches. We \"touch\" every part\n# that the target writer re-serialized by nudging properties back to their\n# own current value (no semantic change), which forces Word to rewrite those\n# elements with the newer attribute ordering.\n\n$d = $word.ActiveDocument\n\n# --- word/document.xml : re-emit the page size / margins (w:sectPr) ---\n$ps = $d.PageSetup\n$ps.PageWidth = $ps.PageWidth\n$ps.PageHeight = $ps.PageHeight\n$ps.TopMargin = $ps.TopMargin\n$ps.RightMargin = $ps.RightMargin\n$ps.BottomMargin = $ps.BottomMargin\n$ps.LeftMargin = $ps.LeftMargin\n$ps.HeaderDistance = $ps.HeaderDistance\n$ps.FooterDistance = $ps.FooterDistance\n$ps.Gutter = $ps.Gutter\n\n# --- word/styles.xml : re-emit docDefaults / latentStyles / every w:style ---\nforeach ($s in $d.Styles) {\n    $s.Priority = $s.Priority\n}\n"}
